# Fix a bug in CalcFullLineExWithMulti: update the per-symbol reel-stop
# statistics table (columns A:F, rows 2-23) to the corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1202, 2, 10, 10, 10, 10)
    3  = @(902, 1, 0, 0, 0, 0)
    4  = @(501, 9, 52, 30, 75, 45)
    5  = @(401, 9, 48, 67, 75, 45)
    6  = @(201, 9, 30, 15, 45, 30)
    7  = @(101, 9, 30, 15, 60, 15)
    8  = @(301, 6, 45, 30, 60, 45)
    9  = @(701, 3, 90, 45, 97, 15)
    10 = @(601, 9, 60, 67, 60, 42)
    12 = @(1201, 2, 10, 10, 10, 10)
    13 = @(1203, 3, 15, 15, 15, 15)
    14 = @(901, 16, 15, 45, 60, 60)
    15 = @(1001, 18, 30, 75, 60, 72)
    16 = @(802, 0, 4, 5, 4, 0)
    17 = @(1, 0, 2, 2, 2, 2)
    18 = @(2, 0, 2, 2, 2, 2)
    20 = @(1101, 0, 15, 30, 30, 0)
    21 = @(502, 0, 4, 0, 0, 0)
    22 = @(402, 0, 0, 4, 0, 0)
    23 = @(602, 0, 0, 4, 0, 9)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
